$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptocurrency price / 1h-volume data (and, for two
# rows, the coin name/link that got re-ordered) to match the latest run of
# the scraper. All data cells in this sheet are stored as literal text
# (matching the source feed's formatting, e.g. "42.882.77" or "0.636"), so
# the D2:E51 range is temporarily switched to a text ("@") number format
# before writing the new values - this stops Excel from "helpfully"
# re-interpreting numeric-looking strings as numbers/dates - and then the
# style is restored to "Normal" afterwards so no stray formatting is left
# behind.

$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '42.882.77'
$ws.Range("E2").Value = '  +4.34%  '
$ws.Range("D3").Value = '2.278.57'
$ws.Range("E3").Value = '  +4.80%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '250.80'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").Value = '0.636'
$ws.Range("E6").Value = '  +3.53%  '
$ws.Range("D7").Value = '71.97'
$ws.Range("E7").Value = '  +9.99%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.658'
$ws.Range("E9").Value = '  +16.59%  '
$ws.Range("D10").Value = '39.09'
$ws.Range("E10").Value = '  +10.82%  '
$ws.Range("D11").Value = '59.78'
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("D12").Value = '0.0971'
$ws.Range("E12").Value = '  +5.07%  '
$ws.Range("D13").Value = '7.48'
$ws.Range("E13").Value = '  +10.06%  '
$ws.Range("D14").Value = '0.104'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").Value = '2.613.46'
$ws.Range("E15").Value = '  +4.67%  '
$ws.Range("D16").Value = '14.91'
$ws.Range("E16").Value = '  +4.30%  '
$ws.Range("D17").Value = '0.888'
$ws.Range("E17").Value = '  +4.77%  '
$ws.Range("D18").Value = '2.276.27'
$ws.Range("E18").Value = '  +4.49%  '
$ws.Range("D19").Value = '42.814.63'
$ws.Range("E19").Value = '  +4.33%  '
$ws.Range("E20").Value = '  +7.62%  '
$ws.Range("D21").Value = '6.32'
$ws.Range("E21").Value = '  +4.06%  '
$ws.Range("D22").Value = '73.15'
$ws.Range("E22").Value = '  +2.38%  '
$ws.Range("D23").Value = '235.95'
$ws.Range("E23").Value = '  +2.85%  '
$ws.Range("E24").Value = '  +5.25%  '
$ws.Range("D25").Value = '4.03'
$ws.Range("E25").Value = '  +5.40%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").Value = '11.40'
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").Value = '2.20'
$ws.Range("E30").Value = '  +5.17%  '
$ws.Range("D31").Value = '167.68'
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").Value = '21.06'
$ws.Range("E32").Value = '  +4.53%  '
$ws.Range("D33").Value = '6.48'
$ws.Range("E33").Value = '  +14.07%  '
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.0806'
$ws.Range("E35").Value = '  +8.87%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = '31.72'
$ws.Range("E36").Value = '  +30.41%  '
$ws.Range("E37").Value = '  +4.17%  '
$ws.Range("D38").Value = '4.53'
$ws.Range("E38").Value = '  +13.62%  '
$ws.Range("D39").Value = '4.78'
$ws.Range("E39").Value = '  +5.98%  '
$ws.Range("D40").Value = '0.0313'
$ws.Range("E40").Value = '  +3.52%  '
$ws.Range("D41").Value = '2.34'
$ws.Range("E41").Value = '  +6.94%  '
$ws.Range("D42").Value = '12.94'
$ws.Range("E42").Value = '  +15.82%  '
$ws.Range("D43").Value = '5.83'
$ws.Range("E43").Value = '  +6.94%  '
$ws.Range("D44").Value = '0.210'
$ws.Range("E44").Value = '  +10.91%  '
$ws.Range("D45").Value = '5.06'
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").Value = '9.21'
$ws.Range("E46").Value = '  +8.26%  '
$ws.Range("D47").Value = '62.30'
$ws.Range("E47").Value = '  +3.64%  '
$ws.Range("E48").Value = '  +2.88%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  +3.50%  '
$ws.Range("E51").Value = '  +4.77%  '

$dataRange.Style = "Normal"
